$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.961.64"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "2.905.65"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "2.902.99"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.428"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "3.385.22"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "61.911.51"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "2.904.58"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.651"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.80%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("E29").Value = "  +6.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").Value = "2.704.81"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "131.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "346.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.15%  "
